# Refresh the coin Price (D) and Volume(1h) (E) columns with the latest
# run's values. The sheet stores these as plain text (e.g. '30.610.57',
# '  +0.42%  '), so for any replacement that LOOKS like a plain number we
# prefix it with a leading apostrophe -- exactly what typing the same text
# into Excel would require -- so the cell stays Text instead of silently
# becoming a Number (which would drop meaningful trailing/leading zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $ws.Range($cellRef).Value = "'" + $text
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

Set-TextValue "D2" '30.623.55'
Set-TextValue "E2" '  +0.48%  '
Set-TextValue "D3" '2.117.40'
Set-TextValue "E3" '  +0.51%  '
Set-TextValue "E4" '  +1.03%  '
Set-TextValue "D5" '340.63'
Set-TextValue "E5" '  +1.86%  '
Set-TextValue "E6" '  +1.07%  '
Set-TextValue "D7" '0.5258'
Set-TextValue "E7" '  +0.05%  '
Set-TextValue "D8" '0.4519'
Set-TextValue "E8" '  -1.36%  '
Set-TextValue "D9" '53.57'
Set-TextValue "E9" '  +0.15%  '
Set-TextValue "D10" '0.09034'
Set-TextValue "E10" '  +0.46%  '
Set-TextValue "D11" '1.171'
Set-TextValue "E11" '  -0.96%  '
Set-TextValue "D12" '24.41'
Set-TextValue "E12" '  +0.14%  '
Set-TextValue "D13" '2.118.04'
Set-TextValue "E13" '  +0.74%  '
Set-TextValue "D14" '6.800'
Set-TextValue "E14" '  +0.17%  '
Set-TextValue "D15" '8.086'
Set-TextValue "E15" '  +2.81%  '
Set-TextValue "D16" '98.24'
Set-TextValue "E16" '  +1.58%  '
Set-TextValue "D17" '0.00001166'
Set-TextValue "E17" '  +3.15%  '
Set-TextValue "D18" '1.014'
Set-TextValue "E18" '  +1.05%  '
Set-TextValue "D19" '0.06710'
Set-TextValue "E19" '  +1.29%  '
Set-TextValue "D20" '19.36'
Set-TextValue "E20" '  -1.03%  '
Set-TextValue "E21" '  +1.02%  '
Set-TextValue "D22" '6.337'
Set-TextValue "E22" '  +0.22%  '
Set-TextValue "D23" '30.692.97'
Set-TextValue "E23" '  +0.54%  '
Set-TextValue "D24" '12.80'
Set-TextValue "E24" '  +3.54%  '
Set-TextValue "D25" '2.386'
Set-TextValue "E25" '  +1.36%  '
Set-TextValue "D26" '2.366.38'
Set-TextValue "E26" '  +0.81%  '
Set-TextValue "D27" '22.40'
Set-TextValue "E27" '  +0.09%  '
Set-TextValue "D28" '165.89'
Set-TextValue "E28" '  +1.16%  '
Set-TextValue "D29" '2.539'
Set-TextValue "E29" '  -1.78%  '
Set-TextValue "D30" '135.75'
Set-TextValue "E30" '  +1.89%  '
Set-TextValue "D31" '1.197'
Set-TextValue "E31" '  -0.02%  '
Set-TextValue "D32" '0.1074'
Set-TextValue "E32" '  -0.14%  '
Set-TextValue "D33" '1.646'
Set-TextValue "E33" '  -1.55%  '
Set-TextValue "D34" '6.383'
Set-TextValue "E34" '  +3.78%  '
Set-TextValue "D35" '3.953'
Set-TextValue "E35" '  +0.62%  '
Set-TextValue "D36" '10.33'
Set-TextValue "E36" '  -1.60%  '
Set-TextValue "D37" '5.890'
Set-TextValue "E37" '  +6.15%  '
Set-TextValue "D38" '0.02652'
Set-TextValue "E38" '  +2.80%  '
Set-TextValue "D39" '0.06840'
Set-TextValue "E39" '  -0.22%  '
Set-TextValue "D40" '0.2325'
Set-TextValue "E40" '  +1.16%  '
Set-TextValue "D41" '12.65'
Set-TextValue "E41" '  -0.92%  '
Set-TextValue "D42" '0.6876'
Set-TextValue "E42" '  -0.44%  '
Set-TextValue "D43" '1.263'
Set-TextValue "E43" '  +1.13%  '
Set-TextValue "D44" '14.83'
Set-TextValue "E44" '  +6.07%  '
Set-TextValue "D45" '0.6429'
Set-TextValue "E45" '  +0.55%  '
Set-TextValue "D46" '2.314'
Set-TextValue "E46" '  -1.63%  '
Set-TextValue "E47" '  +8.75%  '
Set-TextValue "D48" '3.710'
Set-TextValue "E48" '  +1.41%  '
Set-TextValue "D49" '1.255'
Set-TextValue "E49" '  +0.79%  '
Set-TextValue "D50" '82.84'
Set-TextValue "E50" '  -0.78%  '
Set-TextValue "D51" '0.07296'
Set-TextValue "E51" '  +2.91%  '
